$wb = $excel.ActiveWorkbook

# --- Sheet "1D NEW" (sheet2.xml) ---
$wsNew = $wb.Worksheets.Item("1D NEW")
$wsNew.Range("M19").Value = "Brute force solved for 4 detectors."
$wsNew.Range("B15").Select()

# --- Sheet "1D TRAINING" (sheet3.xml) ---
$wsTraining = $wb.Worksheets.Item("1D TRAINING")

$wsTraining.Range("B7").Value = "training_set4"
$wsTraining.Range("C7").Value = 200
$wsTraining.Range("D7").Value = 200
$wsTraining.Range("E7").Value = "0.2-5"
$wsTraining.Range("F7").Value = "0.2-5"
$wsTraining.Range("G7").Value = 60
$wsTraining.Range("H7").Value = 0
$wsTraining.Range("I7").Value = 20
$wsTraining.Range("J7").Value = 30
$wsTraining.Range("K7").Value = 2500

# Column widths (engine quantizes to nearest 1/6 of a character, so we
# pick the ColumnWidth input whose rounded result lands on the nearest
# representable width to the target OOXML width).
$wsTraining.Columns.Item(2).ColumnWidth = 13.0
$wsTraining.Columns.Item(5).ColumnWidth = 14.5
$wsTraining.Columns.Item(6).ColumnWidth = 14.333333333333334
$wsTraining.Columns.Item(7).ColumnWidth = 10.166666666666666
$wsTraining.Columns.Item(8).ColumnWidth = 11.166666666666666
$wsTraining.Columns.Item(9).ColumnWidth = 10.166666666666666
$wsTraining.Columns.Item(10).ColumnWidth = 9.666666666666666

# Activate "1D TRAINING" last so it becomes the active/selected tab
$wsTraining.Activate()
$wsTraining.Range("K8").Select()
